# "first crack at 5/17 stuff"
# Append the 2016-05-17 (serial 42507) row of train-run stats to the
# "Data" sheet. The downstream scatter charts on the "Completion Stats"
# and "Trip Length Stats" chartsheets already reference Data!$A$2:$A$99
# (and the matching B:J columns), so once the sheet gains real values in
# row 34 those series pick the new point up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$row = 34

$ws.Cells.Item($row, 1).Value = 42507                   # A - date (2016-05-17)
$ws.Cells.Item($row, 2).Value = 143                      # B - scheduled
$ws.Cells.Item($row, 3).Value = 135                      # C - completed
$ws.Cells.Item($row, 4).Value = 0                         # D - cancelled (weather, etc.)
$ws.Cells.Item($row, 5).Value = 8                         # E - cancelled (other)
$ws.Cells.Item($row, 6).Value = 135                      # F - completed (total)
$ws.Cells.Item($row, 7).Value = 0.94405594405594406       # G - completion %
$ws.Cells.Item($row, 8).Value = 43.071445221369565         # H - avg metric
$ws.Cells.Item($row, 9).Value = 34.833333335118368         # I - avg metric
$ws.Cells.Item($row, 10).Value = 67.399999997578561        # J - avg metric

# Matches the author leaving the cursor parked on the next empty row
# after typing the new data in.
$ws.Range("A35").Select() | Out-Null
